$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: replace the old "habilidades_linguagens_em" label in C1 with
#     "habilidades" and give it its own bold/white/blue header look ---
$c1 = $ws.Range("C1")
$c1.Value = "habilidades"
$c1.Font.Name = "Calibri"
$c1.Font.Bold = $true
$c1.Font.Size = 14
$c1.Font.Color = 16777215
$c1.Interior.Color = 12611584
$c1.HorizontalAlignment = -4108
$c1.VerticalAlignment = -4108
$c1.WrapText = $true

# Row 1 is now shorter since the header text is short
$ws.Rows.Item(1).RowHeight = 37.5

# Leave the selection on C1, like the author did
$c1.Select() | Out-Null

# Page setup: portrait A4
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
